$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.655.17"
$ws.Range("E2").Value = "  +3.56%  "
$ws.Range("D3").Value = "3.462.06"
$ws.Range("E3").Value = "  +3.94%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "578.63"
$ws.Range("E5").Value = "  +4.70%  "
$ws.Range("D6").Value = "156.64"
$ws.Range("E6").Value = "  +3.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.467.31"
$ws.Range("E8").Value = "  +3.97%  "
$ws.Range("D9").Value = "0.555"
$ws.Range("E9").Value = "  +4.95%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "4.062.75"
$ws.Range("E13").Value = "  +4.37%  "
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("E15").Value = "  +9.19%  "
$ws.Range("D16").Value = "27.89"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "64.675.54"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("D18").Value = "3.455.59"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  +4.61%  "
$ws.Range("D21").Value = "398.09"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "0.549"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").Value = "73.07"
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  +24.40%  "
$ws.Range("D27").Value = "9.46"
$ws.Range("E27").Value = "  +5.85%  "
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "5.98"
$ws.Range("E30").Value = "  +8.27%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "6.74"
$ws.Range("E31").Value = "  +6.26%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("D34").Value = "23.86"
$ws.Range("E34").Value = "  +3.90%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "161.24"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +7.16%  "
$ws.Range("D40").Value = "27.73"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "2.906.46"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").Value = "41.87"
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("D48").Value = "23.39"
$ws.Range("E48").Value = "  +6.00%  "
$ws.Range("D49").Value = "2.18"
$ws.Range("E49").Value = "  +22.46%  "
$ws.Range("D50").Value = "0.864"
$ws.Range("E50").Value = "  +6.06%  "
$ws.Range("E51").Value = "  +3.82%  "
